$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 524.8646210551354
$ws.Range("D2").Value = 123.0273338117717
$ws.Range("F2").Value = 442
$ws.Range("G2").Value = 478
$ws.Range("H2").Value = 570
$ws.Range("C3").Value = 36.74248162610181
$ws.Range("D3").Value = 6.50264959182696
$ws.Range("F3").Value = 31.95
$ws.Range("G3").Value = 36.92
$ws.Range("H3").Value = 41.1
$ws.Range("C4").Value = 2.055739882564128
$ws.Range("D4").Value = 2.579004167187117
$ws.Range("F4").Value = 0.7
$ws.Range("G4").Value = 1.36
$ws.Range("H4").Value = 2.51
$ws.Range("C5").Value = 322.7356937847636
$ws.Range("D5").Value = 9.094335613087898
$ws.Range("F5").Value = 317.69
$ws.Range("G5").Value = 323.21
$ws.Range("H5").Value = 329.32
$ws.Range("C6").Value = 23.17586993857387
$ws.Range("D6").Value = 3.641445351245375
$ws.Range("F6").Value = 20.63
$ws.Range("G6").Value = 22.62
$ws.Range("H6").Value = 25.45
$ws.Range("C7").Value = -76.09269390708312
$ws.Range("D7").Value = 22.57171239534686
$ws.Range("F7").Value = -92
$ws.Range("G7").Value = -72
$ws.Range("H7").Value = -58
$ws.Range("C8").Value = 7.625479137831975
$ws.Range("D8").Value = 6.888120590547293
$ws.Range("F8").Value = 8
$ws.Range("C9").Value = 9.248521444527144
$ws.Range("D9").Value = 1.659694962714828
$ws.Range("C10").Value = 867.8277592894087
$ws.Range("D10").Value = 0.4614728077845632
$ws.Range("C11").Value = 0.5263110209265871
$ws.Range("D11").Value = 0.5714624939927891
$ws.Range("C12").Value = 22.75806666374785
$ws.Range("D12").Value = 12.29785793976647
$ws.Range("C13").Value = 0.6728080750816856
$ws.Range("D13").Value = 0.7505323477745689
$ws.Range("C14").Value = 1.830958331433843
$ws.Range("D14").Value = 1.667744089936605
$ws.Range("C15").Value = 93.35269390708277
$ws.Range("D15").Value = 22.57171239534685
$ws.Range("F15").Value = 75.25999999999999
$ws.Range("G15").Value = 89.25999999999999
$ws.Range("H15").Value = 109.26
$ws.Range("C16").Value = -85.3279545708132
$ws.Range("D16").Value = 20.25890582329711
$ws.Range("F16").Value = -101.6389203414338
$ws.Range("G16").Value = -83.33195619988427
$ws.Range("H16").Value = -69.6389203414338
$ws.Range("C17").Value = -77.70247543298126
$ws.Range("D17").Value = 25.00224010122128
$ws.Range("F17").Value = -92.57382219273629
$ws.Range("G17").Value = -72.46183611348224
$ws.Range("H17").Value = -58.75746206410165
